$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NumberError")

$ws.Range("B2").Value = "Wed Nov 05 11:13:51 EST 2025"
$ws.Range("B3").Value = "Wed Nov 05 11:14:00 EST 2025"
$ws.Range("B4").Value = "Wed Nov 05 11:14:07 EST 2025"
$ws.Range("B5").Value = "Wed Nov 05 11:14:13 EST 2025"
$ws.Range("B6").Value = "Wed Nov 05 11:14:19 EST 2025"
$ws.Range("B7").Value = "Wed Nov 05 11:14:24 EST 2025"
$ws.Range("B8").Value = "Wed Nov 05 11:14:31 EST 2025"
$ws.Range("B9").Value = "Wed Nov 05 11:14:37 EST 2025"
$ws.Range("B10").Value = "Wed Nov 05 11:14:42 EST 2025"
$ws.Range("B11").Value = "Wed Nov 05 11:14:47 EST 2025"
$ws.Range("B12").Value = "Wed Nov 05 11:14:52 EST 2025"
$ws.Range("B13").Value = "Wed Nov 05 11:14:57 EST 2025"
$ws.Range("B14").Value = "Wed Nov 05 11:15:03 EST 2025"
$ws.Range("B15").Value = "Wed Nov 05 11:15:08 EST 2025"
$ws.Range("B16").Value = "Wed Nov 05 11:15:13 EST 2025"
$ws.Range("B17").Value = "Wed Nov 05 11:15:18 EST 2025"
$ws.Range("B18").Value = "Wed Nov 05 11:15:23 EST 2025"
$ws.Range("B19").Value = "Wed Nov 05 11:15:29 EST 2025"
$ws.Range("B20").Value = "Wed Nov 05 11:15:35 EST 2025"
$ws.Range("B21").Value = "Wed Nov 05 11:15:40 EST 2025"
$ws.Range("B22").Value = "Wed Nov 05 11:15:45 EST 2025"
$ws.Range("B23").Value = "Wed Nov 05 11:15:51 EST 2025"
$ws.Range("B24").Value = "Wed Nov 05 11:15:56 EST 2025"
$ws.Range("B25").Value = "Wed Nov 05 11:16:03 EST 2025"
$ws.Range("B26").Value = "Wed Nov 05 11:16:09 EST 2025"
$ws.Range("B27").Value = "Wed Nov 05 11:16:15 EST 2025"
$ws.Range("B28").Value = "Wed Nov 05 11:16:20 EST 2025"
$ws.Range("B29").Value = "Wed Nov 05 11:16:24 EST 2025"
$ws.Range("B30").Value = "Wed Nov 05 11:16:30 EST 2025"
$ws.Range("B31").Value = "Wed Nov 05 11:16:35 EST 2025"
$ws.Range("B32").Value = "Wed Nov 05 11:16:40 EST 2025"
$ws.Range("B33").Value = "Wed Nov 05 11:16:45 EST 2025"
$ws.Range("B34").Value = "Wed Nov 05 11:16:51 EST 2025"
$ws.Range("B35").Value = "Wed Nov 05 11:16:57 EST 2025"
$ws.Range("B36").Value = "Wed Nov 05 11:17:02 EST 2025"
$ws.Range("B37").Value = "Wed Nov 05 11:17:07 EST 2025"
$ws.Range("B38").Value = "Wed Nov 05 11:17:13 EST 2025"
$ws.Range("B39").Value = "Wed Nov 05 11:17:19 EST 2025"
$ws.Range("B40").Value = "Wed Nov 05 11:17:25 EST 2025"
$ws.Range("B44").Value = "Wed Nov 05 11:17:30 EST 2025"
$ws.Range("B45").Value = "Wed Nov 05 11:17:35 EST 2025"
$ws.Range("B46").Value = "Wed Nov 05 11:17:41 EST 2025"
$ws.Range("B47").Value = "Wed Nov 05 11:17:47 EST 2025"
$ws.Range("B48").Value = "Wed Nov 05 11:17:52 EST 2025"
$ws.Range("B49").Value = "Wed Nov 05 11:17:57 EST 2025"
$ws.Range("B50").Value = "Wed Nov 05 11:18:03 EST 2025"
$ws.Range("B51").Value = "Wed Nov 05 11:18:08 EST 2025"
$ws.Range("B52").Value = "Wed Nov 05 11:18:14 EST 2025"
$ws.Range("B53").Value = "Wed Nov 05 11:18:19 EST 2025"
$ws.Range("B54").Value = "Wed Nov 05 11:18:24 EST 2025"
$ws.Range("B55").Value = "Wed Nov 05 11:18:31 EST 2025"
$ws.Range("B56").Value = "Wed Nov 05 11:18:37 EST 2025"
$ws.Range("B57").Value = "Wed Nov 05 11:18:41 EST 2025"
$ws.Range("B58").Value = "Wed Nov 05 11:18:47 EST 2025"
$ws.Range("B62").Value = "Wed Nov 05 11:18:53 EST 2025"
$ws.Range("B63").Value = "Wed Nov 05 11:18:58 EST 2025"
$ws.Range("B64").Value = "Wed Nov 05 11:19:03 EST 2025"
$ws.Range("B65").Value = "Wed Nov 05 11:19:08 EST 2025"
$ws.Range("B66").Value = "Wed Nov 05 11:19:14 EST 2025"
$ws.Range("B67").Value = "Wed Nov 05 11:19:20 EST 2025"
